$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.8222690533928816
$ws.Range("J2").Value = 0.8222690533928814
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.03257366666666667
$ws.Range("N2").Value = 0.097721
$ws.Range("O2").Value = 0.001227793554179957
$ws.Range("P2").Value = 0.001227793554179957
$ws.Range("Q2").Value = 0.07513935987277777
$ws.Range("R2").Value = 0.676254238855
$ws.Range("S2").Value = 0.001009576643557435
$ws.Range("T2").Value = 0.001009576643557435
$ws.Range("I3").Value = 0.8222690533928816
$ws.Range("J3").Value = 0.8222690533928814
$ws.Range("O3").Value = 0.7662385783512358
$ws.Range("P3").Value = 0.7662385783512359
$ws.Range("S3").Value = 0.630054270493978
$ws.Range("T3").Value = 0.6300542704939779
$ws.Range("I4").Value = 0.8222690533928816
$ws.Range("J4").Value = 0.8222690533928814
$ws.Range("M4").Value = 6.169174666666667
$ws.Range("N4").Value = 18.507524
$ws.Range("O4").Value = 0.2325336280945842
$ws.Range("P4").Value = 0.2325336280945842
$ws.Range("Q4").Value = 14.23075394429111
$ws.Range("R4").Value = 128.07678549862
$ws.Range("S4").Value = 0.1912052062553461
$ws.Range("T4").Value = 0.1912052062553461
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.4985973333333333
$ws.Range("H5").Value = 1.495792
$ws.Range("I5").Value = 0.1777309466071186
$ws.Range("J5").Value = 0.1777309466071185
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.03257366666666667
$ws.Range("N5").Value = 0.097721
$ws.Range("O5").Value = 0.001227793554179957
$ws.Range("P5").Value = 0.001227793554179957
$ws.Range("Q5").Value = 0.01624114333688889
$ws.Range("R5").Value = 0.146170290032
$ws.Range("S5").Value = 0.0002182169106225223
$ws.Range("T5").Value = 0.0002182169106225223
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.4985973333333333
$ws.Range("H6").Value = 1.495792
$ws.Range("I6").Value = 0.1777309466071186
$ws.Range("J6").Value = 0.1777309466071185
$ws.Range("O6").Value = 0.7662385783512358
$ws.Range("P6").Value = 0.7662385783512359
$ws.Range("Q6").Value = 10.13573539206933
$ws.Range("R6").Value = 91.22161852862399
$ws.Range("S6").Value = 0.1361843078572579
$ws.Range("T6").Value = 0.1361843078572579
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.4985973333333333
$ws.Range("H7").Value = 1.495792
$ws.Range("I7").Value = 0.1777309466071186
$ws.Range("J7").Value = 0.1777309466071185
$ws.Range("M7").Value = 6.169174666666667
$ws.Range("N7").Value = 18.507524
$ws.Range("O7").Value = 0.2325336280945842
$ws.Range("P7").Value = 0.2325336280945842
$ws.Range("Q7").Value = 3.075934037667555
$ws.Range("R7").Value = 27.683406339008
$ws.Range("S7").Value = 0.04132842183923811
$ws.Range("T7").Value = 0.04132842183923809
